# electricity sensitivity analysis set up
# Renames/restructures the hog-fuel / black-liquor / biosludge rows on the
# Fuels sheet, moving the "moisture fraction" helper cell (column G) down
# from the "dry" row to the corresponding "wet"/derived row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 33: "dry hog fuel" -> "hog fuel (dry)"; drop its G (moisture) value ---
$ws.Range("A33").Value = "hog fuel (dry)"
$ws.Range("G33").ClearContents()

# --- Row 34: "wet hog fuel" -> "hog fuel"; formulas now reference $G34, and
#     the moisture fraction value (0.5) moves here as a new G34 cell ---
$ws.Range("A34").Value = "hog fuel"
$ws.Range("B34").Formula = '=B33*(1-$G34)'
$ws.Range("C34").Formula = '=C33*(1-$G34)'
$ws.Range("D34").Formula = '=D33*(1-$G34)'
$ws.Range("G34").Value = 0.5

# --- Row 35: "dry black liquor" text unchanged; drop its G (moisture) value ---
$ws.Range("A35").Value = "dry black liquor"
$ws.Range("G35").ClearContents()

# --- Row 36: "wet black liquor" -> "strong black liquor"; formulas now
#     reference $G36, and the moisture fraction value (0.2) moves here as a
#     new G36 cell ---
$ws.Range("A36").Value = "strong black liquor"
$ws.Range("B36").Formula = '=B35*(1-$G36)'
$ws.Range("C36").Formula = '=C35*(1-$G36)'
$ws.Range("D36").Formula = '=D35*(1-$G36)'
$ws.Range("E36").Formula = '=G36'
$ws.Range("G36").Value = 0.2

# --- Row 37: "dry biosludge" -> "biosludge (dry)" ---
$ws.Range("A37").Value = "biosludge (dry)"

# --- Row 38: "wet biosludge" -> "biosludge" ---
$ws.Range("A38").Value = "biosludge"

# --- Update the selected/active cell to match the author's saved view ---
$ws.Range("C34").Select()
